$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.735.47"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "3.065.51"
$ws.Range("E3").Value = "  +3.63%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "384.97"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "103.62"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "37.22"
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "0.0867"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "3.550.05"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "18.75"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "7.79"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "3.083.48"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "0.983"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "10.56"
$ws.Range("E18").Value = "  -5.87%  "
$ws.Range("D19").Value = "51.792.59"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "3.16"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "12.56"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "70.29"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "270.18"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "3.18"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").Value = "8.44"
$ws.Range("E26").Value = "  +7.83%  "
$ws.Range("D27").Value = "27.16"
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("E28").Value = "  +5.41%  "
$ws.Range("D29").Value = "7.29"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "34.65"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("D35").Value = "50.49"
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +5.71%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "0.289"
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "17.15"
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("D42").Value = "129.31"
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("D45").Value = "3.74"
$ws.Range("E45").Value = "  +5.54%  "
$ws.Range("D46").Value = "22.21"
$ws.Range("E46").Value = "  +3.63%  "
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").Value = "  +6.39%  "
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "2.052.33"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "3.369.49"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").Value = "0.0321"
$ws.Range("E51").Value = "  +0.56%  "
